$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 10 (Objetivos) B/C content - was pointing at wrong text in the source file
$ws.Range("B10").Value = 'Fornecer os conhecimentos necessários sobre os aspectos fundamentais de Microbiologia e Bioquímica Microbiana e sua importância nos estudos sobre Ecologia dos Microrganismos. Fornecer conhecimentos sobre o papel e utilização dos microrganismos nos processos biológicos de interesse à Engenharia Ambiental.'
$ws.Range("C10").Value = 'Fornecer os conhecimentos necessários sobre os aspectos fundamentais de Microbiologia e Bioquímica Microbiana e sua importância nos estudos sobre Ecologia dos Microrganismos. Fornecer conhecimentos sobre o papel e utilização dos microrganismos nos processos biológicos de interesse à Engenharia Ambiental.'

# Insert a new blank row at 13, shifting old rows 13-23 down to 14-24
$ws.Rows.Item(13).Insert()

# New row 13: move the 'Docentes responsaveis' value into B13/C13
$ws.Range("B13").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("C13").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("A13").Clear()
# Copy B/C column formatting from the (already correctly styled) row below onto row 13
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 14 (was 13): A14 stays 'Programa resumido:'; set the PT short-syllabus content
$ws.Range("B14").Value = 'Diversidade metabólica; cultivo e crescimento microbiano; isolamento microbiano; ecossistemas microbianos; biorremediação e biodeterioração  microbiana; bioindicadores.'
$ws.Range("C14").Value = 'Diversidade metabólica; cultivo e crescimento microbiano; isolamento microbiano; ecossistemas microbianos; biorremediação e biodeterioração  microbiana; bioindicadores.'

# Row 16 (was 15): fix Programa: B/C content (was wrongly showing a date)
$ws.Range("B16").Value = 'Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos,proteínas e ácidos nucleicos.–Diversidade metabólica: Micro-organismos autotróficos e heterotróficos; glicólise; fermentações; respiração; via das pentoses-fosfato; fotossíntese. –Cultivo e crescimento microbiano: Nutrição microbiana; meios de cultura; fatores ambientais; reprodução e crescimento; medidas e controle de crescimento microbiano. –Isolamento microbiano: Técnicas e meios de isolamento.–Ecossistemas microbianos: Diversidade microbiana e ciclos biogeoquímicos. –Biorremediação e biodeterioração microbiana: Lixiviação bacteriana de metais; bioacumulação e biotransformação microbiana de metais; biodegradação de materiais lignocelulósicos; biodegradação de hidrocarbonetos; biodeterioração de monumentos históricos. –Bioindicadores: Bioindicadores de qualidade de água, ar e solo.'
$ws.Range("C16").Value = 'Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos,proteínas e ácidos nucleicos.–Diversidade metabólica: Micro-organismos autotróficos e heterotróficos; glicólise; fermentações; respiração; via das pentoses-fosfato; fotossíntese. –Cultivo e crescimento microbiano: Nutrição microbiana; meios de cultura; fatores ambientais; reprodução e crescimento; medidas e controle de crescimento microbiano. –Isolamento microbiano: Técnicas e meios de isolamento.–Ecossistemas microbianos: Diversidade microbiana e ciclos biogeoquímicos. –Biorremediação e biodeterioração microbiana: Lixiviação bacteriana de metais; bioacumulação e biotransformação microbiana de metais; biodegradação de materiais lignocelulósicos; biodegradação de hidrocarbonetos; biodeterioração de monumentos históricos. –Bioindicadores: Bioindicadores de qualidade de água, ar e solo.'

# Row 19 (was 18): fix Metodo: B/C content (was wrongly showing Objetivos text)
$ws.Range("B19").Value = 'Duas provas escritas (P1 e P2) distribuídas no semestre.'
$ws.Range("C19").Value = 'Duas provas escritas (P1 e P2) distribuídas no semestre.'

# Row 20 (was 19): Criterio: content shifts to the correct text
$ws.Range("B20").Value = 'MF=Média finalMF = (P1 + P2) / 2'
$ws.Range("C20").Value = 'MF=Média finalMF = (P1 + P2) / 2'

# Row 21 (was 20): Norma de recuperacao: content shifts to the correct text
$ws.Range("B21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0.'
$ws.Range("C21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0.'

# Row 22 (was 21): Bibliografia: new bibliography text
$ws.Range("B22").Value = 'Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edição, 2016. - Nelson, D.; Cox, M. Princípios de Bioquímica de Lehninger. Artmed Editora. 6a  Edição, 2014.- Pratt, C.; Cornely, K. Bioquímica essencial. Guanabara Koogan. 1a  Edição, 2006. - Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edição. 2010.- Cooper, G.M. A Célula – Uma Abordagem molecular. Artmed Editora Ltda. 3a  Edição. 2007.- Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Maier, R. Environmental Microbiology. Academic Press. 2000. - Jordening, H.; Winter, J. Environmental Biotechnology. Concepts and Applications. Wiley-VCH. 2005. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edição. 2009.- Tortora, G.; Burdell, B.; Case, C. Microbiology. An Introduction. Pearson Benjamin Cummings. 10a  Edição. 2010.'
$ws.Range("C22").Value = 'Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edição, 2016. - Nelson, D.; Cox, M. Princípios de Bioquímica de Lehninger. Artmed Editora. 6a  Edição, 2014.- Pratt, C.; Cornely, K. Bioquímica essencial. Guanabara Koogan. 1a  Edição, 2006. - Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edição. 2010.- Cooper, G.M. A Célula – Uma Abordagem molecular. Artmed Editora Ltda. 3a  Edição. 2007.- Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Maier, R. Environmental Microbiology. Academic Press. 2000. - Jordening, H.; Winter, J. Environmental Biotechnology. Concepts and Applications. Wiley-VCH. 2005. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edição. 2009.- Tortora, G.; Burdell, B.; Case, C. Microbiology. An Introduction. Pearson Benjamin Cummings. 10a  Edição. 2010.'

# Rows 15, 17, 18, 23, 24 keep the text they inherited from the shift

$ws.Range("A1").Select() | Out-Null
